$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Sector column (K) for the two data rows with "Tech"
$ws.Range("K2").Value = "Tech"
$ws.Range("K3").Value = "Tech"

# Update the active selection to match the edited cell
$ws.Range("K3").Select()
